# Auto-generated script applying scheduled-runner market data refresh
# to the Maduin_Profits leve-profit workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 83.545456
$ws.Range("I33").Value = 83.545456
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 83.545456
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = 145.454544

$ws.Range("H38").Value = 296.6875
$ws.Range("I38").Value = 32.923077
$ws.Range("J38").Value = 1439.6666
$ws.Range("K38").Value = 98.76923099999999
$ws.Range("L38").Value = 4318.9998
$ws.Range("M38").Value = 273.230769
$ws.Range("N38").Value = -5062.9998

$ws.Range("H70").Value = 1923.6666
$ws.Range("I70").Value = 1399.3334
$ws.Range("J70").Value = 2011.0555
$ws.Range("K70").Value = 4198.0002
$ws.Range("L70").Value = 6033.166499999999
$ws.Range("M70").Value = -3928.0002
$ws.Range("N70").Value = -6573.166499999999

$ws.Range("H73").Value = 1923.6666
$ws.Range("I73").Value = 1399.3334
$ws.Range("J73").Value = 2011.0555
$ws.Range("K73").Value = 4198.0002
$ws.Range("L73").Value = 6033.166499999999
$ws.Range("M73").Value = -3262.0002
$ws.Range("N73").Value = -7905.166499999999

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("N86").Value = 0

$ws.Range("H88").Value = 3105.5557
$ws.Range("I88").Value = 2334
$ws.Range("J88").Value = 3491.3333
$ws.Range("K88").Value = 2334
$ws.Range("L88").Value = 3491.3333
$ws.Range("M88").Value = -1928
$ws.Range("N88").Value = -4303.3333

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("N89").Value = 0

$ws.Range("H91").Value = 3105.5557
$ws.Range("I91").Value = 2334
$ws.Range("J91").Value = 3491.3333
$ws.Range("K91").Value = 2334
$ws.Range("L91").Value = 3491.3333
$ws.Range("M91").Value = -930
$ws.Range("N91").Value = -6299.3333

$ws.Range("H135").Value = 1000
$ws.Range("I135").Value = 1000
$ws.Range("K135").Value = 9000
$ws.Range("M135").Value = -6465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").ClearContents()
$ws.Range("N17").Value = 0

$ws.Range("H32").Value = 4414.3335
$ws.Range("I32").Value = 4414.3335
$ws.Range("K32").Value = 4414.3335
$ws.Range("M32").Value = -4127.3335

$ws.Range("H45").Value = 2707.7693
$ws.Range("I45").Value = 967
$ws.Range("J45").Value = 3230
$ws.Range("K45").Value = 967
$ws.Range("L45").Value = 3230
$ws.Range("M45").Value = -590
$ws.Range("N45").Value = -3984

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 950
$ws.Range("I8").Value = 950
$ws.Range("K8").Value = 950
$ws.Range("M8").Value = -810

$ws.Range("H107").Value = 925.8333
$ws.Range("I107").Value = 931.2
$ws.Range("K107").Value = 931.2
$ws.Range("M107").Value = 988.8

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 798.5
$ws.Range("I2").Value = 995
$ws.Range("J2").Value = 602
$ws.Range("K2").Value = 995
$ws.Range("L2").Value = 602
$ws.Range("M2").Value = -882
$ws.Range("N2").Value = -828

$ws.Range("H4").Value = 2436.3635
$ws.Range("J4").Value = 4933.3335
$ws.Range("L4").Value = 4933.3335
$ws.Range("N4").Value = -5157.3335

$ws.Range("H11").Value = 4760.3335
$ws.Range("J11").Value = 2140.5
$ws.Range("L11").Value = 2140.5
$ws.Range("N11").Value = -2420.5

$ws.Range("H86").Value = 11686717
$ws.Range("I86").Value = 17526998
$ws.Range("J86").Value = 6153
$ws.Range("K86").Value = 17526998
$ws.Range("L86").Value = 6153
$ws.Range("M86").Value = -17525875
$ws.Range("N86").Value = -8399

$ws.Range("H89").Value = 11686717
$ws.Range("I89").Value = 17526998
$ws.Range("J89").Value = 6153
$ws.Range("K89").Value = 87634990
$ws.Range("L89").Value = 30765
$ws.Range("M89").Value = -87629374
$ws.Range("N89").Value = -41997

$ws.Range("H99").Value = 4525
$ws.Range("I99").Value = 4525
$ws.Range("K99").Value = 4525
$ws.Range("M99").Value = -3027

$ws.Range("H107").Value = 823.5625
$ws.Range("I107").Value = 848.4666999999999
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 848.4666999999999
$ws.Range("L107").Value = 450
$ws.Range("M107").Value = 1071.5333
$ws.Range("N107").Value = -4290

$ws.Range("H126").Value = 4525
$ws.Range("I126").Value = 4525
$ws.Range("K126").Value = 13575
$ws.Range("M126").Value = -11105

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 150
$ws.Range("I87").Value = 150
$ws.Range("K87").Value = 450
$ws.Range("M87").Value = 798

$ws.Range("H90").Value = 150
$ws.Range("I90").Value = 150
$ws.Range("K90").Value = 1350
$ws.Range("M90").Value = 4890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9981.833000000001
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10312

$ws.Range("H126").Value = 10883.2
$ws.Range("J126").Value = 12138.667
$ws.Range("L126").Value = 36416.001
$ws.Range("N126").Value = -41356.001

$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 7000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H7").Value = 3075.75
$ws.Range("I7").Value = 3467.6667
$ws.Range("J7").Value = 1900
$ws.Range("K7").Value = 3467.6667
$ws.Range("L7").Value = 1900
$ws.Range("M7").Value = -3355.6667
$ws.Range("N7").Value = -2124

$ws.Range("H16").Value = 8500
$ws.Range("I16").Value = 7000
$ws.Range("K16").Value = 7000
$ws.Range("M16").Value = -6830

$ws.Range("H46").Value = 3722.0386
$ws.Range("I46").Value = 3265.1667
$ws.Range("K46").Value = 3265.1667
$ws.Range("M46").Value = -3077.1667

$ws.Range("H126").Value = 3075.75
$ws.Range("I126").Value = 3467.6667
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 10403.0001
$ws.Range("L126").Value = 5700
$ws.Range("M126").Value = -7933.000100000001
$ws.Range("N126").Value = -10640

$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 2707.1667
$ws.Range("I132").Value = 1699.6666
$ws.Range("K132").Value = 5098.9998
$ws.Range("M132").Value = -2568.9998

$ws.Range("H136").Value = 1341.8889
$ws.Range("I136").Value = 680.4
$ws.Range("K136").Value = 2041.2
$ws.Range("M136").Value = 508.8000000000002
